# "Update ILogger of Serilog" — re-key the Yeniemlak.Az "Column / Solution"
# pair: the "+" that used to sit in column Q (Column) actually belongs in
# column R (Solution); clear Q and put the "+" (or "-") into R instead.
# Row 8 additionally gets a brand-new "Qalib" label in column Q.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database Service")

for ($r = 3; $r -le 29; $r++) {
    if ($r -eq 8) {
        $ws.Cells.Item($r, 17).Value = "Qalib"
    } else {
        $ws.Cells.Item($r, 17).Value = $null
    }

    if ($r -eq 14) {
        $ws.Cells.Item($r, 18).Value = "-"
    } else {
        $ws.Cells.Item($r, 18).Value = "+"
    }
}

# Restore the view state recorded in the saved workbook: scrolled so row 8
# is at the top, with U20 as the active selection.
$ws.Activate()
$ws.Range("A8").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("U20").Select()
